# Updates FFXIV leve-profit calculations (currentAveragePrice / Leve price / profit
# columns H-N) per the latest market-board data pull for affected rows across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38: Hi-Potion of Strength
$ws.Range("H38").Value = 457.41666
$ws.Range("I38").Value = 319
$ws.Range("J38").Value = 872.6667
$ws.Range("K38").Value = 957
$ws.Range("L38").Value = 2618.0001
$ws.Range("M38").Value = -585
$ws.Range("N38").Value = -3362.0001
# Row 51: Shark Oil
$ws.Range("H51").Value = 1688.5555
$ws.Range("J51").Value = 1650
$ws.Range("L51").Value = 1650
$ws.Range("N51").Value = -2618
# Row 58: Mega-Potion of Vitality
$ws.Range("H58").Value = 1769.4445
$ws.Range("I58").Value = 1154.1666
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 3462.4998
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -3312.4998
$ws.Range("N58").Value = -9300
# Row 134: Crocodileskin Index
$ws.Range("H134").Value = 84890
$ws.Range("J134").Value = 84890
$ws.Range("L134").Value = 84890
$ws.Range("N134").Value = -95030
# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 1193.6451
$ws.Range("I137").Value = 1061.9048
$ws.Range("J137").Value = 1470.3
$ws.Range("K137").Value = 3185.7144
$ws.Range("L137").Value = 4410.9
$ws.Range("M137").Value = -635.7143999999998
$ws.Range("N137").Value = -9510.9
# Row 138: Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2585.0833
$ws.Range("I138").Value = 3237.3076
$ws.Range("J138").Value = 2404.681
$ws.Range("K138").Value = 9711.9228
$ws.Range("L138").Value = 7214.043
$ws.Range("M138").Value = -4571.9228
$ws.Range("N138").Value = -17494.043

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Range("H32").Value = 393173.4
$ws.Range("I32").Value = 440923.72
$ws.Range("J32").Value = 73246.2
$ws.Range("K32").Value = 440923.72
$ws.Range("L32").Value = 73246.2
$ws.Range("M32").Value = -440636.72
$ws.Range("N32").Value = -73820.2

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Molybdenum Ingot
$ws.Range("H105").Value = 10419291
$ws.Range("I105").Value = 17858956
$ws.Range("J105").Value = 3760
$ws.Range("K105").Value = 17858956
$ws.Range("L105").Value = 3760
$ws.Range("M105").Value = -17857209
$ws.Range("N105").Value = -7254
# Row 134: Ruthenium Ingot
$ws.Range("H134").Value = 2948.4849
$ws.Range("I134").Value = 3225.5789
$ws.Range("J134").Value = 2572.4285
$ws.Range("K134").Value = 9676.736699999999
$ws.Range("L134").Value = 7717.2855
$ws.Range("M134").Value = -7141.736699999999
$ws.Range("N134").Value = -12787.2855
# Row 135: Ruthenium War Axe
$ws.Range("H135").Value = 53072.727
$ws.Range("J135").Value = 47088.89
$ws.Range("L135").Value = 47088.89
$ws.Range("N135").Value = -57228.89
# Row 137: Cobalt Tungsten Khukuri
$ws.Range("H137").Value = 60779
$ws.Range("J137").Value = 60779
$ws.Range("L137").Value = 60779
$ws.Range("N137").Value = -70979

$ws = $wb.Worksheets.Item("CRP")
# Row 28: Iron Lance
$ws.Range("H28").Value = 78333.336
$ws.Range("J28").Value = 78333.336
$ws.Range("L28").Value = 78333.336
$ws.Range("N28").Value = -78823.336
# Row 31: Walnut Lumber
$ws.Range("H31").Value = 3294.1018
$ws.Range("I31").Value = 1011.55884
$ws.Range("J31").Value = 6398.36
$ws.Range("K31").Value = 1011.55884
$ws.Range("L31").Value = 6398.36
$ws.Range("M31").Value = -716.55884
$ws.Range("N31").Value = -6988.36
# Row 34: Walnut Lumber
$ws.Range("H34").Value = 3294.1018
$ws.Range("I34").Value = 1011.55884
$ws.Range("J34").Value = 6398.36
$ws.Range("K34").Value = 1011.55884
$ws.Range("L34").Value = 6398.36
$ws.Range("M34").Value = -809.55884
$ws.Range("N34").Value = -6802.36
# Row 99: Pine Lumber
$ws.Range("H99").Value = 1999.9412
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1999.9412
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1999.9412
$ws.Range("N99").Value = -4995.9412
$ws.Range("M99").ClearContents()
# Row 126: Red Pine Lumber
$ws.Range("H126").Value = 1999.9412
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1999.9412
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 5999.8236
$ws.Range("N126").Value = -10939.8236
$ws.Range("M126").ClearContents()
# Row 134: Ceiba Lumber
$ws.Range("H134").Value = 1305.3334
$ws.Range("I134").Value = 1035.4286
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 3106.2858
$ws.Range("L134").Value = 6750
$ws.Range("M134").Value = -571.2857999999997
$ws.Range("N134").Value = -11820
# Row 135: Ceiba Wings
$ws.Range("H135").Value = 57010
$ws.Range("J135").Value = 57010
$ws.Range("L135").Value = 57010
$ws.Range("N135").Value = -67150

$ws = $wb.Worksheets.Item("CUL")
# Row 5: Maple Syrup
$ws.Range("H5").Value = 893.9032
$ws.Range("I5").Value = 419
$ws.Range("J5").Value = 1311.2424
$ws.Range("K5").Value = 1257
$ws.Range("L5").Value = 3933.7272
$ws.Range("M5").Value = -1145
$ws.Range("N5").Value = -4157.7272
# Row 68: Fermented Butter
$ws.Range("H68").Value = 2308.1633
$ws.Range("I68").Value = 3437.5293
$ws.Range("J68").Value = 1708.1875
$ws.Range("K68").Value = 10312.5879
$ws.Range("L68").Value = 5124.5625
$ws.Range("M68").Value = -9501.5879
$ws.Range("N68").Value = -6746.5625
# Row 71: Fermented Butter
$ws.Range("H71").Value = 2308.1633
$ws.Range("I71").Value = 3437.5293
$ws.Range("J71").Value = 1708.1875
$ws.Range("K71").Value = 30937.7637
$ws.Range("L71").Value = 15373.6875
$ws.Range("M71").Value = -26881.7637
$ws.Range("N71").Value = -23485.6875
# Row 104: Doman Tea
$ws.Range("H104").Value = 3725.4443
$ws.Range("J104").Value = 3725.4443
$ws.Range("L104").Value = 11176.3329
$ws.Range("N104").Value = -16418.3329
# Row 105: Chirashi-zushi
$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -20242
# Row 106: Jerked Jhammel
$ws.Range("H106").Value = 8800
$ws.Range("J106").Value = 8800
$ws.Range("L106").Value = 26400
$ws.Range("N106").Value = -28292
# Row 107: Frantoio Oil
$ws.Range("H107").Value = 2080.4546
$ws.Range("I107").Value = 382.08334
$ws.Range("J107").Value = 2554.4187
$ws.Range("K107").Value = 1146.25002
$ws.Range("L107").Value = 7663.256100000001
$ws.Range("M107").Value = 773.7499800000001
$ws.Range("N107").Value = -11503.2561
# Row 122: Northern Sea Salt
$ws.Range("H122").Value = 7636.2144
$ws.Range("J122").Value = 34099.332
$ws.Range("L122").Value = 306893.988
$ws.Range("N122").Value = -311793.988
# Row 131: Tsai tou Vounou
$ws.Range("H131").Value = 793.8095
$ws.Range("I131").Value = 347
$ws.Range("J131").Value = 1200
$ws.Range("K131").Value = 1041
$ws.Range("L131").Value = 3600
$ws.Range("M131").Value = 3999
$ws.Range("N131").Value = -13680
# Row 135: Royal Maple Syrup
$ws.Range("H135").Value = 893.9032
$ws.Range("I135").Value = 419
$ws.Range("J135").Value = 1311.2424
$ws.Range("K135").Value = 3771
$ws.Range("L135").Value = 11801.1816
$ws.Range("M135").Value = -1236
$ws.Range("N135").Value = -16871.1816

$ws = $wb.Worksheets.Item("GSM")
# Row 132: Lar Ingot
$ws.Range("H132").Value = 3270.8235
$ws.Range("I132").Value = 3034
$ws.Range("J132").Value = 3343.6924
$ws.Range("K132").Value = 9102
$ws.Range("L132").Value = 10031.0772
$ws.Range("M132").Value = -6572
$ws.Range("N132").Value = -15091.0772

$ws = $wb.Worksheets.Item("LTW")
# Row 33: Goatskin Wristbands
$ws.Range("H33").Value = 7258.5
$ws.Range("I33").Value = 6500
$ws.Range("K33").Value = 6500
$ws.Range("M33").Value = -6210
# Row 38: Skull Eyepatch
$ws.Range("H38").Value = 4000
$ws.Range("I38").Value = 3000
$ws.Range("K38").Value = 3000
$ws.Range("M38").Value = -2590
# Row 136: Br'aax Leather
$ws.Range("H136").Value = 3969255.5
$ws.Range("I136").Value = 1079.56
$ws.Range("J136").Value = 9804808
$ws.Range("K136").Value = 3238.68
$ws.Range("L136").Value = 29414424
$ws.Range("M136").Value = -688.6799999999998
$ws.Range("N136").Value = -29419524

$ws = $wb.Worksheets.Item("WVR")
# Row 41: Linen Halfgloves
$ws.Range("H41").Value = 12333.333
$ws.Range("J41").Value = 12333.333
$ws.Range("L41").Value = 12333.333
$ws.Range("N41").Value = -13113.333
# Row 136: Sarcenet Cloth
$ws.Range("H136").Value = 3138.1667
$ws.Range("I136").Value = 2880.3
$ws.Range("J136").Value = 3782.8333
$ws.Range("K136").Value = 8640.900000000001
$ws.Range("L136").Value = 11348.4999
$ws.Range("M136").Value = -6090.900000000001
$ws.Range("N136").Value = -16448.4999
